# The commit changes the content of cell A1 (header) from "cc" to "ID".
# (A2's numeric value of 12345678 is left untouched.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"

# Put the cursor back on the header cell (A1) rather than leaving the
# previously-saved selection on A2.
$ws.Range("A1").Select()
